$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 59; existing rows 59:88 shift down to 60:89
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly record
$ws.Range("A59").Value = 11
$ws.Range("B59").Value = "Vega Monumental Concepción"
$ws.Range("C59").Value = "Bíobío"
$ws.Range("D59").Value = 44460
$ws.Range("E59").Value = 8
$ws.Range("F59").Value = 100112003
$ws.Range("G59").Value = "Ajo"
$ws.Range("H59").Value = "Chino"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 400
$ws.Range("K59").Value = 15000
$ws.Range("L59").Value = 15500
$ws.Range("M59").Value = 15250
$ws.Range("N59").Value = '$/caja 10 kilos'
$ws.Range("O59").Value = "China"
$ws.Range("P59").Value = 1525
$ws.Range("Q59").Value = 10
$ws.Range("R59").Value = "Hortaliza"

Write-Host "Row 59 inserted and populated."
